{"js": "// Replace the date line and the 25 two-digit multiplication problems with\n// their new values, in document order. Each old string is unique in the\n// document, so a plain text search + in-place replace keeps all paragraph /\n// run formatting (fonts, sizes, justification) untouched.\nconst replacements = [\n  [\"2026-01-31 Saturday\", \"2026-02-01 Sunday\"],\n  [\"78\u00d711=\", \"65\u00d777=\"],\n  [\"11\u00d736=\", \"85\u00d796=\"],\n  [\"69\u00d797=\", \"78\u00d771=\"],\n  [\"23\u00d749=\", \"51\u00d753=\"],\n  [\"33\u00d734=\", \"13\u00d775=\"],\n  [\"87\u00d730=\", \"11\u00d715=\"],\n  [\"18\u00d745=\", \"13\u00d743=\"],\n  [\"46\u00d740=\", \"97\u00d747=\"],\n  [\"26\u00d711=\", \"94\u00d785=\"],\n  [\"68\u00d779=\", \"84\u00d797=\"],\n  [\"70\u00d780=\", \"16\u00d790=\"],\n  [\"64\u00d713=\", \"79\u00d739=\"],\n  [\"34\u00d725=\", \"39\u00d758=\"],\n  [\"66\u00d717=\", \"94\u00d736=\"],\n  [\"53\u00d719=\", \"69\u00d769=\"],\n  [\"36\u00d764=\", \"50\u00d739=\"],\n  [\"65\u00d753=\", \"31\u00d725=\"],\n  [\"71\u00d738=\", \"70\u00d797=\"],\n  [\"84\u00d730=\", \"80\u00d767=\"],\n  [\"25\u00d746=\", \"94\u00d721=\"],\n  [\"18\u00d743=\", \"67\u00d751=\"],\n  [\"32\u00d750=\", \"80\u00d721=\"],\n  [\"88\u00d789=\", \"33\u00d738=\"],\n  [\"19\u00d774=\", \"26\u00d726=\"],\n  [\"19\u00d775=\", \"23\u00d765=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 two-digit multiplication problems with\n# their new values, in document order. Each old string is unique in the\n# document, so Find/Replace (wdReplaceAll = 2) on the whole-document range\n# leaves all paragraph / run formatting (fonts, sizes, justification)\n# untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2026-01-31 Saturday'; New = '2026-02-01 Sunday' },\n    @{ Old = '78\u00d711='; New = '65\u00d777=' },\n    @{ Old = '11\u00d736='; New = '85\u00d796=' },\n    @{ Old = '69\u00d797='; New = '78\u00d771=' },\n    @{ Old = '23\u00d749='; New = '51\u00d753=' },\n    @{ Old = '33\u00d734='; New = '13\u00d775=' },\n    @{ Old = '87\u00d730='; New = '11\u00d715=' },\n    @{ Old = '18\u00d745='; New = '13\u00d743=' },\n    @{ Old = '46\u00d740='; New = '97\u00d747=' },\n    @{ Old = '26\u00d711='; New = '94\u00d785=' },\n    @{ Old = '68\u00d779='; New = '84\u00d797=' },\n    @{ Old = '70\u00d780='; New = '16\u00d790=' },\n    @{ Old = '64\u00d713='; New = '79\u00d739=' },\n    @{ Old = '34\u00d725='; New = '39\u00d758=' },\n    @{ Old = '66\u00d717='; New = '94\u00d736=' },\n    @{ Old = '53\u00d719='; New = '69\u00d769=' },\n    @{ Old = '36\u00d764='; New = '50\u00d739=' },\n    @{ Old = '65\u00d753='; New = '31\u00d725=' },\n    @{ Old = '71\u00d738='; New = '70\u00d797=' },\n    @{ Old = '84\u00d730='; New = '80\u00d767=' },\n    @{ Old = '25\u00d746='; New = '94\u00d721=' },\n    @{ Old = '18\u00d743='; New = '67\u00d751=' },\n    @{ Old = '32\u00d750='; New = '80\u00d721=' },\n    @{ Old = '88\u00d789='; New = '33\u00d738=' },\n    @{ Old = '19\u00d774='; New = '26\u00d726=' },\n    @{ Old = '19\u00d775='; New = '23\u00d765=' }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
